$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rows 2-4 rotate (ffffd46d8d01, ffffffddbc0160, 24c6431a)
# and the row that now ends up as row 4 (24c6431a) flips from
# "Handed back: in sync with en-US" to "Ready for handoff" with a new date.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$wsOverview.Range("B2").Value = "e2e\ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$wsOverview.Range("G2").Value = "2016-09-04 07:09:09"

$wsOverview.Range("A3").Value = "ffffffddbc0160-bb20-4166-867a-0e53176b3425.md"
$wsOverview.Range("B3").Value = "e2e\ffffffddbc0160-bb20-4166-867a-0e53176b3425.md"

$wsOverview.Range("A4").Value = "24c6431a-56db-4170-9db2-35732e528a46.md"
$wsOverview.Range("B4").Value = "e2e\24c6431a-56db-4170-9db2-35732e528a46.md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-04 07:11:10"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') { $hl.TextToDisplay = "e2e\ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md" }
    elseif ($addr -eq '$B$3') { $hl.TextToDisplay = "e2e\ffffffddbc0160-bb20-4166-867a-0e53176b3425.md" }
    elseif ($addr -eq '$B$4') { $hl.TextToDisplay = "e2e\24c6431a-56db-4170-9db2-35732e528a46.md" }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$wsZh.Range("G2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-04 07:09:00"
$wsZh.Range("I2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.md"
$wsZh.Range("J2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-04 07:09:28"

$wsZh.Range("A3").Value = "ffffffddbc0160-bb20-4166-867a-0e53176b3425.md"
$wsZh.Range("F3").Value = "True"

$wsZh.Range("A4").Value = "24c6431a-56db-4170-9db2-35732e528a46.md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "24c6431a-56db-4170-9db2-35732e528a46.a599af984170e6d076d7b900d2d85fc4189a7f88.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-04 07:11:01"
$wsZh.Range("I4").Value = "24c6431a-56db-4170-9db2-35732e528a46.md"
$wsZh.Range("J4").Value = "24c6431a-56db-4170-9db2-35732e528a46.a599af984170e6d076d7b900d2d85fc4189a7f88.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-09-04 07:10:35"
$wsZh.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/779c3872e8b254c4e52d8eb21104eb20cb3ea4c3/e2e/24c6431a-56db-4170-9db2-35732e528a46.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/769863a80c4ee30e3393e025b784783a5e4b229f/e2e/24c6431a-56db-4170-9db2-35732e528a46.md."

$wsZh.Columns.Item(16).ColumnWidth = 39.17

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md" }
    elseif ($addr -eq '$I$2') { $hl.TextToDisplay = "230d5acb-504b-4637-9f8e-a5c285449e76.md" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffffddbc0160-bb20-4166-867a-0e53176b3425.md" }
    elseif ($addr -eq '$I$3') { $hl.TextToDisplay = "230d5acb-504b-4637-9f8e-a5c285449e76.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "24c6431a-56db-4170-9db2-35732e528a46.md" }
    elseif ($addr -eq '$I$4') { $hl.TextToDisplay = "24c6431a-56db-4170-9db2-35732e528a46.md" }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md"
$wsDe.Range("G2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-04 07:09:09"
$wsDe.Range("I2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.md"
$wsDe.Range("J2").Value = "230d5acb-504b-4637-9f8e-a5c285449e76.f54858a962c9264f22fab170293849181df9a130.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-04 07:09:36"

$wsDe.Range("A3").Value = "ffffffddbc0160-bb20-4166-867a-0e53176b3425.md"
$wsDe.Range("F3").Value = "True"

$wsDe.Range("A4").Value = "24c6431a-56db-4170-9db2-35732e528a46.md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "24c6431a-56db-4170-9db2-35732e528a46.a599af984170e6d076d7b900d2d85fc4189a7f88.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-04 07:11:10"
$wsDe.Range("I4").Value = "24c6431a-56db-4170-9db2-35732e528a46.md"
$wsDe.Range("J4").Value = "24c6431a-56db-4170-9db2-35732e528a46.a599af984170e6d076d7b900d2d85fc4189a7f88.de-de.xlf"
$wsDe.Range("K4").Value = "2016-09-04 07:10:42"
$wsDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/779c3872e8b254c4e52d8eb21104eb20cb3ea4c3/e2e/24c6431a-56db-4170-9db2-35732e528a46.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/769863a80c4ee30e3393e025b784783a5e4b229f/e2e/24c6431a-56db-4170-9db2-35732e528a46.md."

$wsDe.Columns.Item(16).ColumnWidth = 39.17

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffffd46d8d01-c2cf-4e44-ab34-61cf8c3e76e4.md" }
    elseif ($addr -eq '$I$2') { $hl.TextToDisplay = "230d5acb-504b-4637-9f8e-a5c285449e76.md" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "ffffffddbc0160-bb20-4166-867a-0e53176b3425.md" }
    elseif ($addr -eq '$I$3') { $hl.TextToDisplay = "230d5acb-504b-4637-9f8e-a5c285449e76.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "24c6431a-56db-4170-9db2-35732e528a46.md" }
    elseif ($addr -eq '$I$4') { $hl.TextToDisplay = "24c6431a-56db-4170-9db2-35732e528a46.md" }
}
